$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "73.508.94"
$ws.Range("E2").Value = "  +1.45%  "
$ws.Range("D3").Value = "3.986.56"
$ws.Range("E3").Value = "  -1.44%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'619.45"
$ws.Range("E5").Value = "  +13.90%  "
$ws.Range("D6").Value = "'168.54"
$ws.Range("E6").Value = "  +10.75%  "
$ws.Range("D7").Value = "'0.682"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +0.78%  "
$ws.Range("E10").Value = "  +8.54%  "
$ws.Range("D11").Value = "'56.07"
$ws.Range("E11").Value = "  +3.84%  "
$ws.Range("E12").Value = "  +1.77%  "
$ws.Range("D13").Value = "'11.18"
$ws.Range("E13").Value = "  +2.04%  "
$ws.Range("D14").Value = "4.625.72"
$ws.Range("E14").Value = "  -1.36%  "
$ws.Range("D15").Value = "3.987.79"
$ws.Range("E15").Value = "  -1.34%  "
$ws.Range("D16").Value = "'1.24"
$ws.Range("E16").Value = "  +2.94%  "
$ws.Range("D17").Value = "'14.09"
$ws.Range("E17").Value = "  -1.97%  "
$ws.Range("D18").Value = "'20.50"
$ws.Range("E18").Value = "  -1.06%  "
$ws.Range("D19").Value = "73.316.38"
$ws.Range("E19").Value = "  +1.23%  "
$ws.Range("E20").Value = "  -0.75%  "
$ws.Range("D21").Value = "'440.48"
$ws.Range("E21").Value = "  -1.89%  "
$ws.Range("D22").Value = "'4.89"
$ws.Range("E22").Value = "  +14.01%  "
$ws.Range("D23").Value = "'96.19"
$ws.Range("E23").Value = "  -1.93%  "
$ws.Range("E24").Value = "  -4.37%  "
$ws.Range("D25").Value = "'14.28"
$ws.Range("E25").Value = "  -2.54%  "
$ws.Range("D26").Value = "'4.09"
$ws.Range("E26").Value = "  -3.73%  "
$ws.Range("D27").Value = "'11.08"
$ws.Range("E27").Value = "  -1.96%  "
$ws.Range("D28").Value = "'10.58"
$ws.Range("E28").Value = "  -2.52%  "
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("D30").Value = "'36.22"
$ws.Range("E30").Value = "  -2.84%  "
$ws.Range("D31").Value = "'7.85"
$ws.Range("E31").Value = "  -1.27%  "
$ws.Range("D32").Value = "'13.72"
$ws.Range("E32").Value = "  +0.67%  "
$ws.Range("E33").Value = "  +17.38%  "
$ws.Range("E34").Value = "  -3.36%  "
$ws.Range("D35").Value = "'48.16"
$ws.Range("E35").Value = "  -2.12%  "
$ws.Range("D36").Value = "'70.93"
$ws.Range("E36").Value = "  +6.01%  "
$ws.Range("D37").Value = "'646.58"
$ws.Range("E37").Value = "  -5.05%  "
$ws.Range("E38").Value = "  -4.49%  "
$ws.Range("E39").Value = "  +1.01%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.147"
$ws.Range("E40").Value = "  -1.46%  "
$ws.Range("B41").Value = "Dai"
$ws.Range("C41").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("D43").Value = "'0.0483"
$ws.Range("E43").Value = "  -2.48%  "
$ws.Range("D44").Value = "'10.68"
$ws.Range("E44").Value = "  -4.48%  "
$ws.Range("D45").Value = "'3.22"
$ws.Range("E45").Value = "  -6.26%  "
$ws.Range("E46").Value = "  +34.69%  "
$ws.Range("E47").Value = "  -1.90%  "
$ws.Range("D48").Value = "'0.000297"
$ws.Range("E48").Value = "  +6.28%  "
$ws.Range("D49").Value = "'3.40"
$ws.Range("E49").Value = "  +2.75%  "
$ws.Range("E50").Value = "  -4.45%  "
$ws.Range("D51").Value = "2.830.36"
$ws.Range("E51").Value = "  +3.17%  "
